$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update column F ("想去人数" / want-to-go count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value = 3184
$ws1.Cells.Item(5, 6).Value = 2226
$ws1.Cells.Item(7, 6).Value = 315
$ws1.Cells.Item(8, 6).Value = 1072
$ws1.Cells.Item(10, 6).Value = 248
$ws1.Cells.Item(11, 6).Value = 468
$ws1.Cells.Item(12, 6).Value = 1162
$ws1.Cells.Item(15, 6).Value = 533
$ws1.Cells.Item(16, 6).Value = 7909
$ws1.Cells.Item(19, 6).Value = 219
$ws1.Cells.Item(23, 6).Value = 541
$ws1.Cells.Item(27, 6).Value = 1527
$ws1.Cells.Item(28, 6).Value = 13
$ws1.Cells.Item(30, 6).Value = 1668
$ws1.Cells.Item(34, 6).Value = 47
$ws1.Cells.Item(36, 6).Value = 278
$ws1.Cells.Item(38, 6).Value = 182
$ws1.Cells.Item(39, 6).Value = 349
$ws1.Cells.Item(41, 6).Value = 221

# Sheet "全部类型" (All types) - same events, update column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(5, 6).Value = 3184
$ws4.Cells.Item(7, 6).Value = 2226
$ws4.Cells.Item(9, 6).Value = 315
$ws4.Cells.Item(10, 6).Value = 1072
$ws4.Cells.Item(13, 6).Value = 248
$ws4.Cells.Item(14, 6).Value = 468
$ws4.Cells.Item(15, 6).Value = 1162
$ws4.Cells.Item(18, 6).Value = 533
$ws4.Cells.Item(19, 6).Value = 7909
$ws4.Cells.Item(23, 6).Value = 219
$ws4.Cells.Item(27, 6).Value = 541
$ws4.Cells.Item(31, 6).Value = 1527
$ws4.Cells.Item(32, 6).Value = 13
$ws4.Cells.Item(34, 6).Value = 1668
$ws4.Cells.Item(38, 6).Value = 47
$ws4.Cells.Item(40, 6).Value = 278
$ws4.Cells.Item(42, 6).Value = 182
$ws4.Cells.Item(43, 6).Value = 349
$ws4.Cells.Item(48, 6).Value = 221
